$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = 2.888950333333333
$ws.Range("H2").Value = 8.666850999999999
$ws.Range("I2").Value = 0.01484500611104232
$ws.Range("J2").Value = 0.01484500611104232
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.011299
$ws.Range("N2").Value = 0.033897
$ws.Range("O2").Value = 0.3524439315012944
$ws.Range("P2").Value = 0.3524439315012944
$ws.Range("Q2").Value = 0.03264224981633333
$ws.Range("R2").Value = 0.293780248347
$ws.Range("S2").Value = 0.005232032316936495
$ws.Range("T2").Value = 0.005232032316936495
$ws.Range("G3").Value = 2.888950333333333
$ws.Range("H3").Value = 8.666850999999999
$ws.Range("I3").Value = 0.01484500611104232
$ws.Range("J3").Value = 0.01484500611104232
$ws.Range("O3").Value = 0.01059504871227008
$ws.Range("P3").Value = 0.01059504871227008
$ws.Range("Q3").Value = 0.0009812801298888889
$ws.Range("R3").Value = 0.008831521168999998
$ws.Range("S3").Value = 0.0001572835628804405
$ws.Range("T3").Value = 0.0001572835628804404
$ws.Range("G4").Value = 2.888950333333333
$ws.Range("H4").Value = 8.666850999999999
$ws.Range("I4").Value = 0.01484500611104232
$ws.Range("J4").Value = 0.01484500611104232
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02042033333333334
$ws.Range("N4").Value = 0.061261
$ws.Range("O4").Value = 0.6369610197864354
$ws.Range("P4").Value = 0.6369610197864355
$ws.Range("Q4").Value = 0.05899332879011111
$ws.Range("R4").Value = 0.530939959111
$ws.Range("S4").Value = 0.00945569023122538
$ws.Range("T4").Value = 0.00945569023122538
$ws.Range("I5").Value = 0.06455031362836973
$ws.Range("J5").Value = 0.06455031362836972
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.011299
$ws.Range("N5").Value = 0.033897
$ws.Range("O5").Value = 0.3524439315012944
$ws.Range("P5").Value = 0.3524439315012944
$ws.Range("Q5").Value = 0.1419377969546667
$ws.Range("R5").Value = 1.277440172592
$ws.Range("S5").Value = 0.02275036631482421
$ws.Range("T5").Value = 0.02275036631482421
$ws.Range("I6").Value = 0.06455031362836973
$ws.Range("J6").Value = 0.06455031362836972
$ws.Range("O6").Value = 0.01059504871227008
$ws.Range("P6").Value = 0.01059504871227008
$ws.Range("S6").Value = 0.0006839137172848888
$ws.Range("T6").Value = 0.0006839137172848886
$ws.Range("I7").Value = 0.06455031362836973
$ws.Range("J7").Value = 0.06455031362836972
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.02042033333333334
$ws.Range("N7").Value = 0.061261
$ws.Range("O7").Value = 0.6369610197864354
$ws.Range("P7").Value = 0.6369610197864355
$ws.Range("Q7").Value = 0.2565197916995556
$ws.Range("R7").Value = 2.308678125296
$ws.Range("S7").Value = 0.04111603359626062
$ws.Range("T7").Value = 0.04111603359626062
$ws.Range("G8").Value = 94.40225766666667
$ws.Range("H8").Value = 283.206773
$ws.Range("I8").Value = 0.4850904066394559
$ws.Range("J8").Value = 0.4850904066394559
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.011299
$ws.Range("N8").Value = 0.033897
$ws.Range("O8").Value = 0.3524439315012944
$ws.Range("P8").Value = 0.3524439315012944
$ws.Range("Q8").Value = 1.066651109375667
$ws.Range("R8").Value = 9.599859984380998
$ws.Range("S8").Value = 0.1709671700495714
$ws.Range("T8").Value = 0.1709671700495714
$ws.Range("G9").Value = 94.40225766666667
$ws.Range("H9").Value = 283.206773
$ws.Range("I9").Value = 0.4850904066394559
$ws.Range("J9").Value = 0.4850904066394559
$ws.Range("O9").Value = 0.01059504871227008
$ws.Range("P9").Value = 0.01059504871227008
$ws.Range("Q9").Value = 0.03206530018744445
$ws.Range("R9").Value = 0.288587701687
$ws.Range("S9").Value = 0.005139556488199939
$ws.Range("T9").Value = 0.005139556488199938
$ws.Range("G10").Value = 94.40225766666667
$ws.Range("H10").Value = 283.206773
$ws.Range("I10").Value = 0.4850904066394559
$ws.Range("J10").Value = 0.4850904066394559
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.02042033333333334
$ws.Range("N10").Value = 0.061261
$ws.Range("O10").Value = 0.6369610197864354
$ws.Range("P10").Value = 0.6369610197864355
$ws.Range("Q10").Value = 1.927725568972556
$ws.Range("R10").Value = 17.349530120753
$ws.Range("S10").Value = 0.3089836801016845
$ws.Range("T10").Value = 0.3089836801016845
$ws.Range("G11").Value = 0.421979
$ws.Range("H11").Value = 1.265937
$ws.Range("I11").Value = 0.002168358784660609
$ws.Range("J11").Value = 0.002168358784660609
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.011299
$ws.Range("N11").Value = 0.033897
$ws.Range("O11").Value = 0.3524439315012944
$ws.Range("P11").Value = 0.3524439315012944
$ws.Range("Q11").Value = 0.004767940721
$ws.Range("R11").Value = 0.042911466489
$ws.Range("S11").Value = 0.0007642248949711535
$ws.Range("T11").Value = 0.0007642248949711536
$ws.Range("G12").Value = 0.421979
$ws.Range("H12").Value = 1.265937
$ws.Range("I12").Value = 0.002168358784660609
$ws.Range("J12").Value = 0.002168358784660609
$ws.Range("O12").Value = 0.01059504871227008
$ws.Range("P12").Value = 0.01059504871227008
$ws.Range("Q12").Value = 0.0001433322003333333
$ws.Range("R12").Value = 0.001289989803
$ws.Range("S12").Value = 0.0000229738669491579
$ws.Range("T12").Value = 0.0000229738669491579
$ws.Range("G13").Value = 0.421979
$ws.Range("H13").Value = 1.265937
$ws.Range("I13").Value = 0.002168358784660609
$ws.Range("J13").Value = 0.002168358784660609
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.02042033333333334
$ws.Range("N13").Value = 0.061261
$ws.Range("O13").Value = 0.6369610197864354
$ws.Range("P13").Value = 0.6369610197864355
$ws.Range("Q13").Value = 0.008616951839666669
$ws.Range("R13").Value = 0.07755256655700001
$ws.Range("S13").Value = 0.001381160022740297
$ws.Range("T13").Value = 0.001381160022740297
$ws.Range("G14").Value = 84.33238866666666
$ws.Range("H14").Value = 252.997166
$ws.Range("I14").Value = 0.4333459148364715
$ws.Range("J14").Value = 0.4333459148364714
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.011299
$ws.Range("N14").Value = 0.033897
$ws.Range("O14").Value = 0.3524439315012944
$ws.Range("P14").Value = 0.3524439315012944
$ws.Range("Q14").Value = 0.9528716595446665
$ws.Range("R14").Value = 8.575844935901999
$ws.Range("S14").Value = 0.1527301379249911
$ws.Range("T14").Value = 0.1527301379249911
$ws.Range("G15").Value = 84.33238866666666
$ws.Range("H15").Value = 252.997166
$ws.Range("I15").Value = 0.4333459148364715
$ws.Range("J15").Value = 0.4333459148364714
$ws.Range("O15").Value = 0.01059504871227008
$ws.Range("P15").Value = 0.01059504871227008
$ws.Range("Q15").Value = 0.02864490135044444
$ws.Range("R15").Value = 0.257804112154
$ws.Range("S15").Value = 0.004591321076955659
$ws.Range("T15").Value = 0.004591321076955659
$ws.Range("G16").Value = 84.33238866666666
$ws.Range("H16").Value = 252.997166
$ws.Range("I16").Value = 0.4333459148364715
$ws.Range("J16").Value = 0.4333459148364714
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.02042033333333334
$ws.Range("N16").Value = 0.061261
$ws.Range("O16").Value = 0.6369610197864354
$ws.Range("P16").Value = 0.6369610197864355
$ws.Range("Q16").Value = 1.722095487369556
$ws.Range("R16").Value = 15.498859386326
$ws.Range("S16").Value = 0.2760244558345247
$ws.Range("T16").Value = 0.2760244558345247
